# Rename sheet "stfiulok" -> "jhagjas" and update numeric cell values
# on that sheet per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("stfiulok")
$ws.Name = "jhagjas"

$updates = @{
    "J1" = 38.54852271080017
    "J2" = 39.10971736907959
    "J3" = 39.3399076461792
    "J4" = 39.47608399391174
    "J5" = 39.79212856292725
    "B6" = 2578
    "D6" = 2557
    "E6" = 20
    "F6" = 13
    "G6" = 99.49416342412451
    "H6" = 99.22390376406675
    "I6" = 0.01283547257876313
    "J6" = 40.30218744277954
    "J7" = 40.1474404335022
    "B8" = 2181
    "E8" = 47
    "H8" = 97.8440366972477
    "I8" = 0.0224824355971897
    "J8" = 40.87517118453979
    "B9" = 1768
    "E9" = 13
    "H9" = 99.26428975664969
    "I9" = 0.01135073779795687
    "J9" = 39.79865431785583
    "J10" = 40.04406905174255
    "J11" = 40.00422048568726
    "B12" = 2539
    "E12" = 1
    "H12" = 99.96059889676911
    "I12" = 0.0003940110323089047
    "J12" = 39.75757908821106
    "B13" = 1806
    "E13" = 12
    "H13" = 99.33518005540166
    "I13" = 0.007242339832869081
    "J13" = 39.88439130783081
    "J14" = 39.70295524597168
    "J15" = 39.24021124839783
    "J16" = 40.54644918441772
    "B17" = 1534
    "E17" = 0
    "H17" = 100
    "I17" = 0
    "J17" = 39.79909944534302
    "J18" = 40.52464532852173
    "J19" = 42.02813577651978
    "J20" = 39.72353029251099
    "J21" = 40.10713529586792
    "B22" = 1519
    "E22" = 2
    "H22" = 99.86824769433466
    "I22" = 0.001318391562294001
    "J22" = 39.75511240959167
    "J23" = 39.36995577812195
    "J24" = 40.043301820755
    "B25" = 1952
    "E25" = 2
    "H25" = 99.89748846745259
    "I25" = 0.007641365257259297
    "J25" = 39.58996868133545
    "J26" = 39.12427878379822
    "B27" = 2877
    "D27" = 2864
    "F27" = 114
    "G27" = 96.1719274680994
    "H27" = 99.58275382475661
    "I27" = 0.04229607250755287
    "J27" = 40.15285205841064
    "J28" = 39.28722906112671
    "B29" = 2064
    "D29" = 1852
    "E29" = 211
    "F29" = 6
    "G29" = 99.67707212055974
    "H29" = 89.77217644207465
    "I29" = 0.1167294244217321
    "J29" = 39.38474154472351
    "B30" = 2944
    "E30" = 6
    "H30" = 99.79612640163099
    "I30" = 0.007111412123264477
    "J30" = 40.32072949409485
    "J31" = 40.25624775886536
    "B32" = 2624
    "D32" = 2621
    "F32" = 27
    "G32" = 98.98036253776435
    "H32" = 99.9237514296607
    "I32" = 0.01094752736881842
    "J32" = 39.3378758430481
    "J33" = 41.54549169540405
    "B34" = 3246
    "E34" = 2
    "H34" = 99.93836671802774
    "I34" = 0.001847290640394089
    "J34" = 43.36691403388977
    "J35" = 39.70415186882019
    "B36" = 3359
    "D36" = 3358
    "F36" = 3
    "G36" = 99.91074085093722
    "I36" = 0.000892325996430696
    "J36" = 40.79715967178345
    "B37" = 3282
    "E37" = 1078
    "H37" = 67.14416336482779
    "I37" = 0.490036231884058
    "J37" = 40.4207935333252
    "J38" = 39.61178946495056
    "J39" = 40.35430788993835
    "J40" = 39.30700588226318
    "B41" = 2485
    "D41" = 2479
    "F41" = 2
    "G41" = 99.91938734381299
    "H41" = 99.7987117552335
    "I41" = 0.00282030620467365
    "J41" = 39.79197692871094
    "J42" = 40.09335041046143
    "B43" = 2052
    "E43" = 7
    "H43" = 99.65870307167235
    "I43" = 0.00682261208576998
    "J43" = 39.88975071907043
    "J44" = 40.35441946983337
    "J45" = 39.83396649360657
    "J46" = 39.27898764610291
    "J47" = 41.00643682479858
    "J48" = 40.17040872573853
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
